# "update sprintu 3, uzupelnienie in/out"
# Fill in the input/output description cells for the ECG_BASELINE module row
# (row 10/11, merged visually with the team-member rows) on the IN_OUT sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D10 = "Dane wejściowe" (input) for ECG_BASELINE
$ws.Range("D10").Value = "Vector <double> od modułu I/O - nieprzefiltrowany sygnał"

# E10 = "Dane wyjściowe" (output) for ECG_BASELINE
$ws.Range("E10").Value = "Vector <double>  - przefiltrowany sygnał z usuniętą linią bazową"

# Leave the view as the author left it: scrolled/selected around the area just edited.
$ws.Activate()
$ws.Range("E11").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
